$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'274.36"
$ws.Range("D3").Value = "'22.96"
$ws.Range("D4").Value = "'6.365"
$ws.Range("D5").Value = "'0.06234"
$ws.Range("D6").Value = "'3.660"
$ws.Range("D7").Value = "'6.679"
$ws.Range("D8").Value = "'1.371"
$ws.Range("D9").Value = "'0.8314"
$ws.Range("D10").Value = "'0.01378"
$ws.Range("D11").Value = "'0.1631"
$ws.Range("D12").Value = "'0.08293"
$ws.Range("D13").Value = "'0.03394"
$ws.Range("D14").Value = "'0.03111"
$ws.Range("B15").Value = "ProBitToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D15").Value = "'0.1269"
$ws.Range("E15").Value = "14ProBitTokenPROB"
$ws.Range("B16").Value = "BitMartToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D16").Value = "'0.09319"
$ws.Range("E16").Value = "15BitMartTokenBMX"
$ws.Range("B17").Value = "MCDex"
$ws.Range("C17").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D17").Value = "'3.891"
$ws.Range("E17").Value = "16MCDexMCB"
$ws.Range("B18").Value = "BitForexToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D18").Value = "'0.001660"
$ws.Range("E18").Value = "17BitForexTokenBF"
$ws.Range("B19").Value = "CoinExToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D19").Value = "'0.04775"
$ws.Range("E19").Value = "18CoinExTokenCET"
$ws.Range("B20").Value = "TigerCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D20").Value = "'0.006384"
$ws.Range("E20").Value = "19TigerCashTCH"
$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D21").Value = "'0.005559"
$ws.Range("E21").Value = "20HotbitTokenHTB"
$ws.Range("B22").Value = "BitKan"
$ws.Range("C22").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D22").Value = "'0.001090"
$ws.Range("E22").Value = "21BitKanKAN"
$ws.Range("B23").Value = "NitroEx"
$ws.Range("C23").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D23").Value = "'0.0001500"
$ws.Range("E23").Value = "22NitroExNTX"
$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").Value = "'3.727"
$ws.Range("E24").Value = "23LEOLEO"
$ws.Range("B25").Value = "BTSEToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D25").Value = "'2.370"
$ws.Range("E25").Value = "24BTSETokenBTSE"
$ws.Range("B26").Value = "BitpandaEcosystemToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D26").Value = "'0.3380"
$ws.Range("E26").Value = "25BitpandaEcosystemTokenBEST"
$ws.Range("D27").Value = "'0.0002681"
$ws.Range("D40").Value = "'0.04690"
$ws.Range("D41").Value = "'0.007020"
$ws.Range("D43").Value = "'0.003330"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"
$ws.Range("D45").Value = "'0.00006257"
$ws.Range("D47").Value = "'0.9005"
$ws.Range("D48").Value = "'0.03022"
$ws.Range("D49").Value = "'0.00002301"
$ws.Range("D50").Value = "'0.01241"
